# Restore the "min <= hour" threshold for rule R30 (sheet "Rules", row 10,
# column C) from 18 back to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
